# The commit replaces the word "for" with "with" in the Professional
# Certificates entry: "Generative AI for Large Language Models (2023)"
# -> "Generative AI with Large Language Models (2023)".
#
# (The rest of the upstream XML diff is just Word's own proofing-mark
# clean-up / run re-splitting noise with no effect on the visible text,
# so it is not something to reproduce through the object model.)

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "Generative AI for Large Language Models",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Generative AI with Large Language Models",
    2
)
